$d = $word.ActiveDocument

$pkgHead = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">'
$pkgTail = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 3: "Configure the build path to include Tomcat and the JDBC connector file" ---
$inner = '<w:body><w:p w14:paraId="10841212" w14:textId="1DEB2FB8" w:rsidR="009A0B97" w:rsidRDefault="009A0B97" w:rsidP="009A0B97">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Configure the build path to include Tomcat and the JDBC connector </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>file</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p></w:body>'
$pkg = $pkgHead + $inner + $pkgTail
$rng = $d.Paragraphs.Item(3).Range
$rng.InsertXML($pkg)

# --- Paragraph 4: "Create basic HTML file to include in the WebINF folder" ---
$inner = '<w:body><w:p w14:paraId="6292642F" w14:textId="755E30A1" w:rsidR="009A0B97" w:rsidRDefault="009A0B97" w:rsidP="009A0B97">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Create basic HTML file to include in the </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>WebINF</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>folder</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p></w:body>'
$pkg = $pkgHead + $inner + $pkgTail
$rng = $d.Paragraphs.Item(4).Range
$rng.InsertXML($pkg)

# --- Paragraph 6: "From the original Class Diagram, the Login and Registration classes will be Servlets" ---
$inner = '<w:body><w:p w14:paraId="77AB6BA6" w14:textId="39D2C693" w:rsidR="004B2F24" w:rsidRDefault="00520872" w:rsidP="004B2F24">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">From the original Class Diagram, the Login and Registration classes will be </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Servlets</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p></w:body>'
$pkg = $pkgHead + $inner + $pkgTail
$rng = $d.Paragraphs.Item(6).Range
$rng.InsertXML($pkg)

# --- Paragraph 7: "Create the database" ---
$inner = '<w:body><w:p w14:paraId="6618C2A7" w14:textId="4B03A393" w:rsidR="00793B5A" w:rsidRDefault="00793B5A" w:rsidP="004B2F24">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Create the </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>database</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p></w:body>'
$pkg = $pkgHead + $inner + $pkgTail
$rng = $d.Paragraphs.Item(7).Range
$rng.InsertXML($pkg)

# --- Paragraph 8: "Make it manually in MySQL Workbench" (only 2nd run " Workbench" changes) ---
$inner = '<w:body><w:p w14:paraId="65A70D00" w14:textId="62AF68DB" w:rsidR="004E13DA" w:rsidRDefault="004E13DA" w:rsidP="004E13DA">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Make it manually in MySQL</w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Workbench</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p></w:body>'
$pkg = $pkgHead + $inner + $pkgTail
$rng = $d.Paragraphs.Item(8).Range
$rng.InsertXML($pkg)

# --- Paragraph 9: "In order to collaborate, we are using Github (through the native Eclipse features)" ---
$inner = '<w:body><w:p w14:paraId="7E569496" w14:textId="375406F4" w:rsidR="00F234FC" w:rsidRDefault="00F234FC" w:rsidP="00F234FC">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>In order to</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> collaborate, we are </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">using </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Github</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (through the native Eclipse features)</w:t></w:r>' +
  '</w:p></w:body>'
$pkg = $pkgHead + $inner + $pkgTail
$rng = $d.Paragraphs.Item(9).Range
$rng.InsertXML($pkg)

# --- Paragraph 10: "Learning how to get Git to work in Eclipse was difficult for all of us" ---
$inner = '<w:body><w:p w14:paraId="2D9B8FA4" w14:textId="37AE558A" w:rsidR="00365313" w:rsidRDefault="008A3D71" w:rsidP="00365313">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Learning how to get Git to work in Eclipse was difficult for all of </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>us</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p></w:body>'
$pkg = $pkgHead + $inner + $pkgTail
$rng = $d.Paragraphs.Item(10).Range
$rng.InsertXML($pkg)

# --- Delete paragraphs 11-17 (the whole "Started by coding together ... As a workaround" block) ---
$delStart = $d.Paragraphs.Item(11).Range.Start
$delEnd = $d.Paragraphs.Item(17).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Text = ""
